$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "35.251.98"
$ws.Range("E2").Value = "  -0.33%  "
$ws.Range("D3").Value = "1.909.44"
$ws.Range("E3").Value = "  -0.27%  "
$ws.Range("E4").Value = "  +0.25%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.720"
$ws.Range("E5").Value = "  +9.44%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "253.29"
$ws.Range("E6").Value = "  +3.58%  "
$ws.Range("E7").Value = "  +0.27%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "40.55"
$ws.Range("E8").Value = "  -1.87%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.359"
$ws.Range("E9").Value = "  +2.25%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "52.30"
$ws.Range("E10").Value = "  -1.40%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0768"
$ws.Range("E11").Value = "  +7.14%  "
$ws.Range("E12").Value = "  -0.61%  "
$ws.Range("D13").Value = "2.191.43"
$ws.Range("E13").Value = "  -0.09%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "12.69"
$ws.Range("E14").Value = "  +4.96%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.717"
$ws.Range("E15").Value = "  +2.30%  "
$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D16").Value = "1.942.33"
$ws.Range("E16").Value = "  +1.43%  "
$ws.Range("B17").Value = "Polkadot"
$ws.Range("C17").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "4.90"
$ws.Range("E17").Value = "  +0.48%  "
$ws.Range("D18").Value = "35.274.29"
$ws.Range("E18").Value = "  -0.21%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "74.12"
$ws.Range("E19").Value = "  +2.89%  "
$ws.Range("D20").Value = "0.0₃0850"
$ws.Range("E20").Value = "  +3.24%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "243.31"
$ws.Range("E21").Value = "  +1.55%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "12.97"
$ws.Range("E22").Value = "  +3.58%  "
$ws.Range("E23").Value = "  +5.07%  "
$ws.Range("E24").Value = "  +0.29%  "
$ws.Range("B25").Value = "Toncoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.36"
$ws.Range("E25").Value = "  +2.77%  "
$ws.Range("B26").Value = "PancakeSwap"
$ws.Range("C26").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.42"
$ws.Range("E26").Value = "  +3.87%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "167.09"
$ws.Range("E27").Value = "  -1.90%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.60"
$ws.Range("E28").Value = "  +1.96%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "18.69"
$ws.Range("E29").Value = "  +1.25%  "
$ws.Range("E30").Value = "  +4.02%  "
$ws.Range("D31").Value = "4.128.76"
$ws.Range("E31").Value = "  +19.46%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.33"
$ws.Range("E32").Value = "  +4.33%  "
$ws.Range("E33").Value = "  +13.85%  "
$ws.Range("B34").Value = "Hedera"
$ws.Range("C34").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0582"
$ws.Range("E34").Value = "  +2.67%  "
$ws.Range("B35").Value = "TrustWalletToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.61"
$ws.Range("E35").Value = "  +21.41%  "
$ws.Range("E36").Value = "  +1.58%  "
$ws.Range("E37").Value = "  -0.64%  "
$ws.Range("E38").Value = "  -1.96%  "
$ws.Range("E39").Value = "  -1.03%  "
$ws.Range("B40").Value = "VeChain"
$ws.Range("C40").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0218"
$ws.Range("E40").Value = "  +4.16%  "
$ws.Range("B41").Value = "InjectiveProtocol"
$ws.Range("C41").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "17.14"
$ws.Range("E41").Value = "  +4.36%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "96.76"
$ws.Range("E42").Value = "  +6.91%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.12"
$ws.Range("E43").Value = "  -0.06%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0642"
$ws.Range("E44").Value = "  -2.88%  "
$ws.Range("D45").Value = "1.336.97"
$ws.Range("E45").Value = "  -0.46%  "
$ws.Range("E46").Value = "  +1.95%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.43"
$ws.Range("E47").Value = "  +0.59%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "6.79"
$ws.Range("E48").Value = "  +2.71%  "
$ws.Range("E49").Value = "  -1.07%  "
$ws.Range("B50").Value = "MultiversX"
$ws.Range("C50").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "45.23"
$ws.Range("E50").Value = "  -6.03%  "
$ws.Range("B51").Value = "Gas"
$ws.Range("C51").Value = "https://coinranking.com/coin/hfw0nnnLtSFc7+gas-gas"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "12.07"
$ws.Range("E51").Value = "  +22.08%  "
